# Cities_by_Region.xlsx — move the Pennsylvania cities that were filed under
# "northeast" (Pittsburgh, Erie, Johnstown, Altoona, New Castle) into the
# "midwest" sheet, inserted in their correctly sorted (HOLC Neighborhoods
# descending) position, and remove them from "northeast".

$wb = $excel.ActiveWorkbook

$midwest   = $wb.Worksheets.Item("midwest")
$northeast = $wb.Worksheets.Item("northeast")

# --- 1. Insert the five rows into "midwest" at their correctly sorted
#        positions (ascending target row order so each insert lands where
#        intended without needing to re-compute offsets). ---

$midwest.Rows.Item(6).Insert()
$midwest.Range("A6").Value = "Pittsburgh, PA"
$midwest.Range("B6").Value = "Pittsburgh"
$midwest.Range("C6").Value = 114

$midwest.Rows.Item(14).Insert()
$midwest.Range("A14").Value = "Erie, PA"
$midwest.Range("B14").Value = "Erie"
$midwest.Range("C14").Value = 47

$midwest.Rows.Item(23).Insert()
$midwest.Range("A23").Value = "Johnstown, PA"
$midwest.Range("B23").Value = "Johnstown"
$midwest.Range("C23").Value = 35

$midwest.Rows.Item(31).Insert()
$midwest.Range("A31").Value = "Altoona, PA"
$midwest.Range("B31").Value = "Altoona"
$midwest.Range("C31").Value = 30

$midwest.Rows.Item(34).Insert()
$midwest.Range("A34").Value = "New Castle, PA"
$midwest.Range("B34").Value = "New Castle"
$midwest.Range("C34").Value = 26

# --- 2. Remove the same five cities from "northeast" (descending row order
#        so earlier deletes don't shift the rows still to be removed). ---

$northeast.Rows.Item(33).Delete()
$northeast.Rows.Item(26).Delete()
$northeast.Rows.Item(24).Delete()
$northeast.Rows.Item(16).Delete()
$northeast.Rows.Item(6).Delete()
